$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "J25" = -10.57102434548657
    "K25" = 3.574792626865744
    "I26" = -10.51602434548657
    "J26" = 3.629792626865739
    "H27" = -10.51602434548659
    "I27" = 3.629792626865722
    "G28" = -10.51602434548657
    "H28" = 3.629792626865739
    "F29" = -10.51602434548659
    "G29" = 3.689792626865739
    "H29" = 5.716794673020033
    "I29" = -2.932891992481572
    "J29" = -5.907106464233138
    "K29" = 2.050717317831413
    "E30" = -10.56702434548656
    "F30" = 3.629792626865722
    "G30" = 5.89679467302004
    "H30" = -2.792891992481557
    "I30" = -5.867106464233117
    "J30" = 2.090717317831433
    "D31" = -7.516024345486572
    "E31" = 5.729792626865716
    "F31" = 7.756794673020053
    "G31" = -1.632891992481575
    "H31" = -5.867603096389431
    "I31" = 2.090222366053069
    "C32" = -2.61602434548657
    "D32" = 6.329792626865739
    "E32" = 6.156794673020033
    "F32" = -2.232891992481555
    "G32" = -5.80710646423312
    "H32" = 2.15071731783143
    "B33" = -3.56602434548657
    "C33" = 5.769792626865737
    "D33" = 6.056794673020022
    "E33" = -2.332891992481564
    "F33" = -5.807106464233115
    "G33" = 2.150717317831436
    "H33" = 0.5746947653272656
    "I33" = -0.3780328061469476
    "J33" = -0.5979284922632784
    "K33" = 1.443932912058642
    "B34" = 2.529792626865742
    "C34" = 5.356794673020033
    "D34" = -0.1328919924815612
    "E34" = -5.307106464233129
    "F34" = 2.39071731783143
    "G34" = 0.5746947653272656
    "H34" = -0.3780328061469476
    "I34" = -0.5979284922632784
    "J34" = 1.443932912058642
    "B35" = 2.956794673020027
    "C35" = 0.06710800751844204
    "D35" = -1.607106464233111
    "E35" = 3.550717317831427
    "F35" = 0.804913431933457
    "G35" = -0.2518180252671272
    "H35" = -0.580456746937102
    "I35" = 1.430045938313753
    "B36" = 0.2479822279162387
    "C36" = 0.4928935357668678
    "D36" = 4.150717317831436
    "E36" = 0.7346947653272622
    "F36" = -0.3180328061469595
    "G36" = -0.5635338037639739
    "H36" = 1.390442373121999
    "B37" = 0.1928935357668848
    "C37" = 4.350717317831439
    "D37" = 1.034694765327259
    "E37" = -0.1180328061469567
    "F37" = -0.4379284922632818
    "G37" = 1.50393291205863
    "H37" = 1.429861316022425
    "I37" = 0.1526454189713746
    "J37" = 0.4656763841019966
    "K37" = -0.03673872975578271
    "B38" = 0.6507173178314358
    "C38" = -1.265305234672738
    "D38" = 3.681967193853055
    "E38" = 1.412071507736727
    "F38" = 1.523932912058654
    "G38" = 1.579861316022431
    "H38" = 0.3026454189713803
    "I38" = 0.5156763841019938
    "J38" = 0.0132612702442145
    "B39" = 0.5556414999948345
    "C39" = 0.7719671938530577
    "D39" = 0.1360272089500689
    "E39" = 2.103932912058638
    "F39" = 2.199861316022421
    "G39" = 0.9826454189713727
    "H39" = 0.9456763841019864
    "I39" = 0.1832612702442162
    "B40" = 0.7219671938530607
    "C40" = 0.162071507736727
    "D40" = 1.853932912058638
    "E40" = 1.929861316022411
    "F40" = 0.5526454189713661
    "G40" = 0.8656763841019881
    "H40" = 0.113261270244223
    "B41" = -1.186928492263277
    "C41" = 0.3639329120586434
    "D41" = -0.1101386839775809
    "E41" = -0.4473545810286339
    "F41" = 0.5156763841019938
    "G41" = 0.2632612702442003
    "H41" = 1.014894517653033
    "I41" = 0.8219948986110523
    "J41" = 0.3112534500261574
    "K41" = 0.5166541826233555
    "B42" = 0.5039362906915836
    "C42" = 0.2797613160224302
    "D42" = -0.3773545810286265
    "E42" = 0.805676384102
    "F42" = 0.3032612702442208
    "G42" = 0.9148945176530248
    "H42" = 0.5219948986110409
    "I42" = 0.01125345002614597
    "J42" = 0.1566541826233561
    "B43" = 0.4898613160224272
    "C43" = -0.1073545810286305
    "D43" = 0.6656763841019995
    "E43" = 0.3632612702442088
    "F43" = 0.9148945176530248
    "G43" = 0.5519948986110705
    "H43" = -0.04874654997384201
    "I43" = 0.116654182623364
    "B44" = -0.298354581028633
    "C44" = 0.8146763841020004
    "D44" = 0.3632612702442088
    "E44" = 1.063894517653054
    "F44" = 0.8709948986110732
    "G44" = 0.2602534500261697
    "H44" = 0.3656541826233735
    "B45" = 0.4656763841019966
    "C45" = 0.2632612702442003
    "D45" = 0.8148945176530447
    "E45" = 0.8219948986110523
    "F45" = 0.3262534500261722
    "G45" = 0.416654182623347
    "H45" = 0.01016824066003608
    "I45" = 0.2452723979283746
    "B46" = -0.116738729755781
    "C46" = 0.7148945176530219
    "D46" = 0.6219948986110637
    "E46" = 0.3112534500261574
    "F46" = 0.5166541826233555
    "G46" = 0.01016824066003608
    "H46" = 0.3452723979283974
    "B47" = 0.5148945176530333
    "C47" = 0.6219948986110637
    "D47" = 0.3112534500261574
    "E47" = 0.5166541826233555
    "F47" = 0.01016824066003608
    "G47" = 0.3452723979283974
    "B48" = 0.221994898611058
    "C48" = 0.1602534500261754
    "D48" = 0.3856541826233695
    "E48" = 0.06116814066005816
    "F48" = 0.505272397928394
    "B49" = -0.1887465499738426
    "C49" = 0.2166541826233441
    "D49" = 0.01016824066003608
    "E49" = 0.545272397928386
    "B50" = 0.2366557702529377
    "C50" = -0.1898317593399668
    "D50" = 0.2052723979283826
    "B51" = -0.3898340444052479
    "C51" = 0.1852723979283866
    "B52" = 0.1452758398526868
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$deletes = @(
    "J45",
    "I46",
    "H47",
    "G48",
    "F49",
    "E50",
    "D51",
    "C52",
    "B53",
)

foreach ($addr in $deletes) {
    $ws.Range($addr).ClearContents()
}
